# Add three "bad word" entries that were missing from the list, which
# pushes all subsequent rows down (dimension grows from C201 to C204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the old row 106 ("chmo") for the two
# missing "buvini ami[ga ske]" entries.
$ws.Rows.Item(106).Insert()
$ws.Rows.Item(106).Insert()

$ws.Range("A106").Value = 9
$ws.Range("B106").Value = "buvini ami"
$ws.Range("C106").Value = 100

$ws.Range("A107").Value = 10
$ws.Range("B107").Value = "buvini amiga ske"
$ws.Range("C107").Value = 100

# Insert one more row before what is now row 116 (old "fuck" entry, now
# shifted to 118) for the missing "foxisha qanchiq" entry.
$ws.Rows.Item(116).Insert()

$ws.Range("A116").Value = 8
$ws.Range("B116").Value = "foxisha qanchiq"
$ws.Range("C116").Value = 79
